$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Data changes: fill in "actual time" (column D) for newly completed topics ---
# Row 9  - dimentionality reduction
$ws.Range("D9").Value = 0.5
# Row 10 - clustering
$ws.Range("D10").Value = 2.5
# Row 13 - course (computer vision)
$ws.Range("D13").Value = 0.5
# Row 14 - cv programming basics
$ws.Range("D14").Value = 0.5

# D5 (=SUM(D6:D11)) and D12 (=SUM(D13:D20)) are formulas and recalculate
# automatically to 9.5 and 1 respectively once the above values are set.

# --- View state: move selection/scroll position down to around row 12-15 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 12
$win.ScrollColumn = 1
$ws.Range("D15").Select()
